$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 308.33334
$ws.Range("I2").Value = 270
$ws.Range("K2").Value = 270
$ws.Range("M2").Value = -157

$ws.Range("H53").Value = 794.46155
$ws.Range("I53").Value = 853.2222
$ws.Range("J53").Value = 662.25
$ws.Range("K53").Value = 853.2222
$ws.Range("L53").Value = 662.25
$ws.Range("M53").Value = -216.2222
$ws.Range("N53").Value = -1936.25

$ws.Range("H103").Value = 6299.5
$ws.Range("J103").Value = 7285
$ws.Range("L103").Value = 21855
$ws.Range("N103").Value = -23027

$ws.Range("H116").Value = 6286.4287
$ws.Range("J116").Value = 5800
$ws.Range("L116").Value = 5800
$ws.Range("N116").Value = -12684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11352.406
$ws.Range("I32").Value = 10388.896
$ws.Range("K32").Value = 10388.896
$ws.Range("M32").Value = -10101.896

$ws.Range("H45").Value = 2719.3125
$ws.Range("I45").Value = 1563.625
$ws.Range("J45").Value = 3875
$ws.Range("K45").Value = 1563.625
$ws.Range("L45").Value = 3875
$ws.Range("M45").Value = -1186.625
$ws.Range("N45").Value = -4629

$ws.Range("H61").Value = 5159.6
$ws.Range("I61").Value = 800
$ws.Range("K61").Value = 800
$ws.Range("M61").Value = -588

$ws.Range("H97").Value = 222
$ws.Range("I97").Value = 222
$ws.Range("K97").Value = 222
$ws.Range("M97").Value = 274

$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

$ws.Range("H132").Value = 2160.111
$ws.Range("I132").Value = 1829.8572
$ws.Range("J132").Value = 3316
$ws.Range("K132").Value = 5489.571599999999
$ws.Range("L132").Value = 9948
$ws.Range("M132").Value = -2959.571599999999
$ws.Range("N132").Value = -15008

$ws.Range("H136").Value = 5159.6
$ws.Range("I136").Value = 800
$ws.Range("K136").Value = 2400
$ws.Range("M136").Value = 150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 758.36365
$ws.Range("I80").Value = 487.5
$ws.Range("K80").Value = 487.5
$ws.Range("M80").Value = 510.5

$ws.Range("H83").Value = 758.36365
$ws.Range("I83").Value = 487.5
$ws.Range("K83").Value = 2437.5
$ws.Range("M83").Value = 2554.5

$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 500
$ws.Range("M86").Value = 623

$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 2500
$ws.Range("M89").Value = 3116

$ws.Range("H99").Value = 1635.7273
$ws.Range("I99").Value = 1499.5
$ws.Range("K99").Value = 1499.5
$ws.Range("M99").Value = -1.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3578.889
$ws.Range("I31").Value = 3387.4285
$ws.Range("K31").Value = 3387.4285
$ws.Range("M31").Value = -3092.4285

$ws.Range("H34").Value = 3578.889
$ws.Range("I34").Value = 3387.4285
$ws.Range("K34").Value = 3387.4285
$ws.Range("M34").Value = -3185.4285

$ws.Range("H98").Value = 56666.668
$ws.Range("J98").Value = 70000
$ws.Range("L98").Value = 70000
$ws.Range("N98").Value = -74492

$ws.Range("H107").Value = 1367.75
$ws.Range("I107").Value = 929
$ws.Range("J107").Value = 2099
$ws.Range("K107").Value = 929
$ws.Range("L107").Value = 2099
$ws.Range("M107").Value = 991
$ws.Range("N107").Value = -5939

$ws.Range("H134").Value = 3695.6538
$ws.Range("I134").Value = 3300.8
$ws.Range("J134").Value = 5011.8335
$ws.Range("K134").Value = 9902.400000000001
$ws.Range("L134").Value = 15035.5005
$ws.Range("M134").Value = -7367.400000000001
$ws.Range("N134").Value = -20105.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.333332
$ws.Range("I2").Value = 16.285715
$ws.Range("K2").Value = 97.71429000000001
$ws.Range("M2").Value = 15.28570999999999

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9094792
$ws.Range("I70").Value = 14289242
$ws.Range("J70").Value = 4503.5
$ws.Range("K70").Value = 14289242
$ws.Range("L70").Value = 4503.5
$ws.Range("M70").Value = -14288972
$ws.Range("N70").Value = -5043.5

$ws.Range("H73").Value = 9094792
$ws.Range("I73").Value = 14289242
$ws.Range("J73").Value = 4503.5
$ws.Range("K73").Value = 14289242
$ws.Range("L73").Value = 4503.5
$ws.Range("M73").Value = -14288306
$ws.Range("N73").Value = -6375.5

$ws.Range("H122").Value = 7464.6
$ws.Range("I122").Value = 8521.286
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 25563.858
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -23113.858
$ws.Range("N122").Value = -19897

$ws.Range("H132").Value = 3243.0952
$ws.Range("I132").Value = 2739.2307
$ws.Range("J132").Value = 4061.875
$ws.Range("K132").Value = 8217.6921
$ws.Range("L132").Value = 12185.625
$ws.Range("M132").Value = -5687.6921
$ws.Range("N132").Value = -17245.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1270.4546
$ws.Range("I22").Value = 897.5
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 897.5
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -602.5
$ws.Range("N22").Value = -5590

$ws.Range("H27").Value = 1270.4546
$ws.Range("I27").Value = 897.5
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 897.5
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -790.5
$ws.Range("N27").Value = -5214

$ws.Range("H46").Value = 1441.5714
$ws.Range("J46").Value = 1499.6
$ws.Range("L46").Value = 1499.6
$ws.Range("N46").Value = -1875.6

$ws.Range("H106").Value = 23190.2
$ws.Range("J106").Value = 23190.2
$ws.Range("L106").Value = 23190.2
$ws.Range("N106").Value = -25714.2

$ws.Range("H132").Value = 4359
$ws.Range("I132").Value = 2102.6
$ws.Range("K132").Value = 6307.799999999999
$ws.Range("M132").Value = -3777.799999999999

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 23332.666
$ws.Range("J63").Value = 23332.666
$ws.Range("L63").Value = 23332.666
$ws.Range("N63").Value = -24580.666

$ws.Range("H66").Value = 23332.666
$ws.Range("J66").Value = 23332.666
$ws.Range("L66").Value = 69997.99800000001
$ws.Range("N66").Value = -76237.99800000001

$ws.Range("H103").Value = 32427.25
$ws.Range("J103").Value = 32427.25
$ws.Range("L103").Value = 32427.25
$ws.Range("N103").Value = -34771.25
